{"js": "// Update the endorsement date from \"March 28, 2023\" to \"July 04, 2023\",\n// and update the travel-report month range from\n// \"January 1 to March 31, 2022.\" to \"April 1 to June 30, 2023.\"\n// Both strings occur twice in the document (the letter is duplicated twice\n// on the page), so every match returned by search() is updated.\n\nconst body = context.document.body;\n\nconst dateResults = body.search(\"March 28, 2023\", { matchCase: true });\ndateResults.load(\"text\");\n\nconst rangeResults = body.search(\"January 1 to March 31, 2022.\", { matchCase: true });\nrangeResults.load(\"text\");\n\nawait context.sync();\n\nfor (let i = 0; i < dateResults.items.length; i++) {\n  dateResults.items[i].insertText(\"July 04, 2023\", Word.InsertLocation.replace);\n}\n\nfor (let i = 0; i < rangeResults.items.length; i++) {\n  rangeResults.items[i].insertText(\"April 1 to June 30, 2023.\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the endorsement date from \"March 28, 2023\" to \"July 04, 2023\",\n# and update the travel-report month range from\n# \"January 1 to March 31, 2022.\" to \"April 1 to June 30, 2023.\"\n# Both strings occur twice in the document (the letter is duplicated twice\n# on the page), so wdReplaceAll (2) updates every occurrence.\n\n$d = $word.ActiveDocument\n\n$find1 = $d.Content.Find\n$find1.Execute(\"March 28, 2023\", $false, $false, $false, $false, $false, $true, 1, $false, \"July 04, 2023\", 2)\n\n$find2 = $d.Content.Find\n$find2.Execute(\"January 1 to March 31, 2022.\", $false, $false, $false, $false, $false, $true, 1, $false, \"April 1 to June 30, 2023.\", 2)\n"}
